# Calculos para UVLO OVLO - adds a new OVLO/UVLO resistor-divider + Zener
# calculation block (rows 16-28, new Table "Tabla3"), updates RB value for
# the 3.3V column in the first calc block, and widens columns C/D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update existing data: RB (D5) changes from 10000 to 5900 for the
#    3.3V column. Downstream formulas in D6/D7 recalculate automatically.
# ---------------------------------------------------------------------
$ws.Range("D5").Value = 5900

# ---------------------------------------------------------------------
# 2) New block: resistor divider for the OVLO sense pin (rows 16-21)
# ---------------------------------------------------------------------
$ws.Range("B16").Value = "Variable"
$ws.Range("C16").Value = "VALOR 5V"
$ws.Range("D16").Value = "VALOR 3,3V"

$ws.Range("B17").Value = "V_O(V)"
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3.3

$ws.Range("B18").Value = "V_f(V)"
$ws.Range("C18").Formula = "=C17*C21/(C20+C21)"
$ws.Range("D18").Formula = "=D17*D21/(D20+D21)"

$ws.Range("B19").Value = "SELA(V)"
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 0

$ws.Range("B20").Value = "R_UP(ohm)"
$ws.Range("C20").Value = 6980
$ws.Range("D20").Value = 6980

$ws.Range("B21").Value = "R_DOWN(ohm)"
$ws.Range("C21").Value = 2210
$ws.Range("D21").Formula = "=2210+1820"

# Apply the same cell styles used by the other two calc tables: the
# "VARIABLE" column uses the 60% - Enfasis2 cell style, and the value
# columns use the Calculo cell style. The last computed resistor value
# (C21) is emphasized with bold+underline, same as the other blocks.
$ws.Range("B17:B21").Style = "60% - Énfasis2"
$ws.Range("C17:D20").Style = "Cálculo"
$ws.Range("D21").Style = "Cálculo"
$ws.Range("C21").Style = "Cálculo"
$ws.Range("C21").Font.Bold = $true
$ws.Range("C21").Font.Underline = $true
$ws.Range("C21").Value = 2210

# Turn B16:D21 into an actual Excel Table ("Tabla3"), matching the style
# used by the other tables on the sheet.
$null = $ws.ListObjects.Add(1, $ws.Range("B16:D21"), 0, 1)
$ws.ListObjects.Item(4).Name = "Tabla3"
$ws.ListObjects.Item("Tabla3").TableStyle = "TableStyleLight8"

# ---------------------------------------------------------------------
# 3) New block: Zener / reference based OVLO_DOWN - OVLO_UP calculation
#    (rows 23-28), not part of a table.
# ---------------------------------------------------------------------
$ws.Range("B23").Value = "V_Z"
$ws.Range("C23").Value = 2.5
$ws.Range("D23").Value = 2.5

$ws.Range("B24").Value = "V_REF"
$ws.Range("C24").Value = 1.25
$ws.Range("D24").Value = 1.22

$ws.Range("B25").Value = "R_UP"
$ws.Range("C25").Value = 10000
$ws.Range("D25").Value = 10000

$ws.Range("B26").Value = "R_DOWN"
$ws.Range("C26").Formula = "=C25*((C23/C24)-1)"
$ws.Range("D26").Formula = "=D25*((D23/D24)-1)"

$ws.Range("C27").Value = "OVLO_DOWN"
$ws.Range("D27").Value = "OVLO_UP"

$ws.Range("B28").Value = "V_OUT"

# ---------------------------------------------------------------------
# 4) Column widths for the new wider labels/values in columns C and D.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 11.0220
$ws.Columns.Item(4).ColumnWidth = 12.5924

# ---------------------------------------------------------------------
# 5) Leave the view focused on the new block, like in the final workbook.
# ---------------------------------------------------------------------
$ws.Range("B28").Select()
